$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_1")

$ws.Range("C3").Value = 9936.588112500001
$ws.Range("D3").Value = 134907.7581017592
$ws.Range("E3").Value = 3514088.075655271
$ws.Range("F3").Value = 1.590786763863847
$ws.Range("K3").Value = 3025.885112537143
$ws.Range("N3").Value = 774.839307540194
$ws.Range("P3").Value = 33.36810169136437
$ws.Range("R3").Value = 0.1670631120265793
$ws.Range("U3").Value = 1177.239333887838
$ws.Range("Y3").Value = 18.09164673579438
